$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 9959
    "F3" = 219
    "F4" = 43
    "F5" = 595
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cell in $updates.Keys) {
        $ws.Range($cell).Value = $updates[$cell]
    }
}
